$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update marking row total count (B11: 3 -> 5)
$ws.Range("B11").Value = 5

# Update total marks (B12: 45 -> 75)
$ws.Range("B12").Value = 75

# Update Corr/Total marks text (E12: "40/84" -> "75/140")
$ws.Range("E12").Value = "75/140"
